$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns B:E (rows 2-51) to Text so numeric-looking strings (e.g. "21.30", "1.000")
# are stored verbatim as text, matching the source data export; cleared back to
# the default (General) style afterwards so no residual formatting remains on the cells.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '27.159.81'
$ws.Range("E2").Value = '  +1.11%  '

$ws.Range("D3").Value = '1.896.91'
$ws.Range("E3").Value = '  +1.93%  '

$ws.Range("D4").Value = '0.9999'
$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").Value = '307.34'
$ws.Range("E5").Value = '  +0.94%  '

$ws.Range("D6").Value = '0.9998'
$ws.Range("E6").Value = '  +0.01%  '

$ws.Range("D7").Value = '0.5156'
$ws.Range("E7").Value = '  +2.22%  '

$ws.Range("D8").Value = '0.3764'
$ws.Range("E8").Value = '  +3.54%  '

$ws.Range("D9").Value = '0.07219'
$ws.Range("E9").Value = '  +0.47%  '

$ws.Range("D10").Value = '21.30'
$ws.Range("E10").Value = '  +2.67%  '

$ws.Range("D11").Value = '0.9065'
$ws.Range("E11").Value = '  +1.31%  '

$ws.Range("D12").Value = '0.07650'
$ws.Range("E12").Value = '  +2.16%  '

$ws.Range("D13").Value = '1.886.01'
$ws.Range("E13").Value = '  +1.20%  '

$ws.Range("D14").Value = '95.22'
$ws.Range("E14").Value = '  +3.03%  '

$ws.Range("D15").Value = '5.270'
$ws.Range("E15").Value = '  +0.83%  '

$ws.Range("D16").Value = '1.000'
$ws.Range("E16").Value = '  +0.04%  '

$ws.Range("D17").Value = '0.000008509'
$ws.Range("E17").Value = '  +0.30%  '

$ws.Range("D18").Value = '14.47'
$ws.Range("E18").Value = '  +2.24%  '

$ws.Range("D19").Value = '0.9997'
$ws.Range("E19").Value = '  -0.03%  '

$ws.Range("D20").Value = '27.178.95'
$ws.Range("E20").Value = '  +1.07%  '

$ws.Range("D21").Value = '5.078'
$ws.Range("E21").Value = '  +0.89%  '

$ws.Range("D22").Value = '2.126.82'
$ws.Range("E22").Value = '  +1.11%  '

$ws.Range("D23").Value = '10.61'
$ws.Range("E23").Value = '  +2.50%  '

$ws.Range("D24").Value = '6.418'
$ws.Range("E24").Value = '  +0.18%  '

$ws.Range("D25").Value = '2.306'
$ws.Range("E25").Value = '  +11.38%  '

$ws.Range("D26").Value = '146.47'
$ws.Range("E26").Value = '  -0.55%  '

$ws.Range("D27").Value = '1.773'
$ws.Range("E27").Value = '  -0.93%  '

$ws.Range("D28").Value = '18.08'
$ws.Range("E28").Value = '  +1.22%  '

$ws.Range("D29").Value = '114.70'
$ws.Range("E29").Value = '  +1.40%  '

$ws.Range("D30").Value = '4.947'
$ws.Range("E30").Value = '  +5.74%  '

$ws.Range("D31").Value = '4.832'
$ws.Range("E31").Value = '  +3.31%  '

$ws.Range("D32").Value = '0.09202'
$ws.Range("E32").Value = '  -0.42%  '

$ws.Range("D33").Value = '0.05087'
$ws.Range("E33").Value = '  -0.15%  '

$ws.Range("D34").Value = '1.253'
$ws.Range("E34").Value = '  +9.10%  '

$ws.Range("D35").Value = '0.7832'
$ws.Range("E35").Value = '  +4.71%  '

$ws.Range("D36").Value = '2.971'
$ws.Range("E36").Value = '  -0.56%  '

$ws.Range("E37").Value = '  +0.24%  '

$ws.Range("D38").Value = '2.630'
$ws.Range("E38").Value = '  +4.54%  '

$ws.Range("B39").Value = 'TheSandbox'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D39").Value = '0.5629'
$ws.Range("E39").Value = '  +2.14%  '

$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").Value = '0.02004'
$ws.Range("E40").Value = '  +0.29%  '

$ws.Range("D41").Value = '1.075'
$ws.Range("E41").Value = '  +0.84%  '

$ws.Range("D42").Value = '9.112'
$ws.Range("E42").Value = '  +6.92%  '

$ws.Range("D43").Value = '6.687'
$ws.Range("E43").Value = '  +2.39%  '

$ws.Range("D44").Value = '118.08'
$ws.Range("E44").Value = '  -0.12%  '

$ws.Range("D45").Value = '0.1510'
$ws.Range("E45").Value = '  +3.04%  '

$ws.Range("D46").Value = '0.4820'
$ws.Range("E46").Value = '  +3.18%  '

$ws.Range("D47").Value = '10.23'
$ws.Range("E47").Value = '  +1.67%  '

$ws.Range("D48").Value = '0.9999'
$ws.Range("E48").Value = '  +0.05%  '

$ws.Range("D49").Value = '1.602'
$ws.Range("E49").Value = '  +2.51%  '

$ws.Range("D50").Value = '37.69'
$ws.Range("E50").Value = '  +2.31%  '

$ws.Range("D51").Value = '64.26'
$ws.Range("E51").Value = '  +1.89%  '

$dataRange.ClearFormats()
